# Adding Vavg and Vpp unit-conversion rows to the "Elec" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: single-quote char, used to build the quote-prefixed Latex formula
# strings (they are literal text beginning with a double-quote, but Excel's
# "quote prefix" cell style is applied the same way the existing rows use it
# - by feeding a value that starts with an apostrophe, which Excel strips
# and replaces with the quote-prefix flag).
$q = [char]39

# Row 125: Voltage [V] -> Voltage [Vp]
$ws.Range("A125").Value = 124
$ws.Range("B125").Value = "ELEC"
$ws.Range("C125").Value = "Voltage [V]"
$ws.Range("D125").Value = "[V]"
$ws.Range("E125").Value = "Voltage [Vp]"
$ws.Range("F125").Value = "[Vp]"
$ws.Range("I125").Value = "sympy.Pow(2,0.5) * P[0]"
$ws.Range("J125").Value = ($q + '"r' + $q + '$Voltage [Vpeak] = \sqrt{2} \times Voltage[V]$' + $q + '"')
$ws.Range("K125").Value = "V"
$ws.Range("L125").Value = 0.001
$ws.Range("M125").Value = 340

# Row 126: Voltage [Vp] -> Voltage [V]
$ws.Range("A126").Value = 125
$ws.Range("B126").Value = "ELEC"
$ws.Range("C126").Value = "Voltage [Vp]"
$ws.Range("D126").Value = "[Vp]"
$ws.Range("E126").Value = "Voltage [V]"
$ws.Range("F126").Value = "[V]"
$ws.Range("I126").Value = " P[0]/sympy.Pow(2,0.5)"
$ws.Range("J126").Value = ($q + '"r' + $q + '$Voltage [V] = \frac{Voltage[Vpeak]}{\sqrt{2}}$' + $q + '"')
$ws.Range("K126").Value = "Vp"
$ws.Range("L126").Value = 0.001
$ws.Range("M126").Value = 340

# Row 127: Voltage [V] -> Voltage [Vpp]
$ws.Range("A127").Value = 125
$ws.Range("B127").Value = "ELEC"
$ws.Range("C127").Value = "Voltage [V]"
$ws.Range("D127").Value = "[V]"
$ws.Range("E127").Value = "Voltage [Vpp]"
$ws.Range("F127").Value = "[Vpp]"
$ws.Range("I127").Value = "2*sympy.Pow(2,0.5) * P[0]"
$ws.Range("J127").Value = ($q + '"r' + $q + '$Voltage [Vpeak-peak] =2 \sqrt{2} \times Voltage[V]$' + $q + '"')
$ws.Range("K127").Value = "V"
$ws.Range("L127").Value = 0.001
$ws.Range("M127").Value = 340

# Row 128: Voltage [Vpp] -> Voltage [V]
$ws.Range("A128").Value = 125
$ws.Range("B128").Value = "ELEC"
$ws.Range("C128").Value = "Voltage [Vpp]"
$ws.Range("D128").Value = "[Vpp]"
$ws.Range("E128").Value = "Voltage [V]"
$ws.Range("F128").Value = "[V]"
$ws.Range("I128").Value = " P[0]/(2*sympy.Pow(2,0.5))"
$ws.Range("J128").Value = ($q + '"r' + $q + '$Voltage [V] = \frac{Voltage[Vpeak-peak]}{2\sqrt{2}}$' + $q + '"')
$ws.Range("K128").Value = "Vpp"
$ws.Range("L128").Value = 0.001
$ws.Range("M128").Value = 340

# Row 129: Voltage [Vavg] -> Voltage [V]
$ws.Range("A129").Value = 125
$ws.Range("B129").Value = "ELEC"
$ws.Range("C129").Value = "Voltage [Vavg]"
$ws.Range("D129").Value = "[Vavg]"
$ws.Range("E129").Value = "Voltage [V]"
$ws.Range("F129").Value = "[V]"
$ws.Range("I129").Value = "sympy.pi*P[0]/(2*sympy.Pow(2,0.5))"
$ws.Range("J129").Value = ($q + '"r' + $q + '$Voltage [V] = \frac{\pi \times Voltage[Vavg]}{2\sqrt{2}}$' + $q + '"')
$ws.Range("K129").Value = "Vavg"
$ws.Range("L129").Value = 0.001
$ws.Range("M129").Value = 340

# Row 130: Voltage [V] -> Voltage [Vavg]
$ws.Range("A130").Value = 125
$ws.Range("B130").Value = "ELEC"
$ws.Range("C130").Value = "Voltage [V]"
$ws.Range("D130").Value = "[V]"
$ws.Range("E130").Value = "Voltage [Vavg]"
$ws.Range("F130").Value = "[Vavg]"
$ws.Range("I130").Value = "(2*sympy.Pow(2,0.5)*P[0])/sympy.pi"
$ws.Range("J130").Value = ($q + '"r' + $q + '$Voltage [Vavg] = \frac{2\sqrt{2}}{\pi} \times Voltage [V]$' + $q + '"')
$ws.Range("K130").Value = "V"
$ws.Range("L130").Value = 0.001
$ws.Range("M130").Value = 340

# Leave the selection where the author ended up after entering the new rows.
$ws.Range("I136").Select()
